$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$r2A = @"
Pipeline(steps=[('scaler', None),
                ('selector',
                 SelectFromModel(estimator=ExtraTreesClassifier(random_state=42))),
                ('model',
                 BaggingClassifier(estimator=MLPClassifier(activation='tanh',
                                                           alpha=1e-05,
                                                           hidden_layer_sizes=(10,),
                                                           learning_rate_init=0.01,
                                                           max_iter=1000,
                                                           random_state=42,
                                                           solver='lbfgs'),
                                   n_estimators=50, random_state=42))])
"@
$ws.Range("A2").Value = $r2A
$r2C = @"
{'selector': SelectFromModel(estimator=ExtraTreesClassifier(random_state=42)), 'scaler': None, 'model__n_estimators': 50, 'model__estimator__solver': 'lbfgs', 'model__estimator__learning_rate_init': 0.01, 'model__estimator__hidden_layer_sizes': (10,), 'model__estimator__alpha': 1e-05, 'model__estimator__activation': 'tanh'}
"@
$ws.Range("C2").Value = $r2C
$r2E = @"
[1 1 0 0 1 0 0 0 0 1 0 1]
"@
$ws.Range("E2").Value = $r2E
$r2F = @"
[0 1 1 1 0 1 1 1 1 1 1 1]
"@
$ws.Range("F2").Value = $r2F
$ws.Range("B2").Value = 0.6857142857142857
$ws.Range("D2").Value = 0.4
$ws.Range("G2").Value = 77
$ws.Range("H2").Value = 0.756904761904762
$ws.Range("I2").Value = 0.02377006313463532
$ws.Range("J2").Value = 0.5822619047619046
$ws.Range("K2").Value = 0.05812372528453319

# Row 3
$r3A = @"
Pipeline(steps=[('scaler', None), ('selector', None),
                ('model',
                 BaggingClassifier(estimator=MLPClassifier(activation='tanh',
                                                           alpha=1e-05,
                                                           hidden_layer_sizes=(5,
                                                                               10,
                                                                               5),
                                                           learning_rate_init=1e-05,
                                                           max_iter=1000,
                                                           random_state=42,
                                                           solver='sgd'),
                                   random_state=42))])
"@
$ws.Range("A3").Value = $r3A
$r3C = @"
{'selector': None, 'scaler': None, 'model__n_estimators': 10, 'model__estimator__solver': 'sgd', 'model__estimator__learning_rate_init': 1e-05, 'model__estimator__hidden_layer_sizes': (5, 10, 5), 'model__estimator__alpha': 1e-05, 'model__estimator__activation': 'tanh'}
"@
$ws.Range("C3").Value = $r3C
$r3E = @"
[1 1 0 1 0 0 1 0 1 1 1 0]
"@
$ws.Range("E3").Value = $r3E
$r3F = @"
[1 1 1 1 1 1 1 1 1 1 1 1]
"@
$ws.Range("F3").Value = $r3F
$ws.Range("B3").Value = 0.6
$ws.Range("D3").Value = 0.7368421052631579
$ws.Range("G3").Value = 69
$ws.Range("H3").Value = 0.7639843279083784
$ws.Range("I3").Value = 0.02699819094802026
$ws.Range("J3").Value = 0.5355033152501507
$ws.Range("K3").Value = 0.0838706479962243

# Row 4
$r4A = @"
Pipeline(steps=[('scaler', None),
                ('selector',
                 <__main__.NamedFeatureSelector object at 0x7f3a6c5ddc10>),
                ('model',
                 BaggingClassifier(estimator=MLPClassifier(hidden_layer_sizes=(10,
                                                                               10,
                                                                               10),
                                                           learning_rate_init=0.01,
                                                           max_iter=1000,
                                                           random_state=42,
                                                           solver='lbfgs'),
                                   random_state=42))])
"@
$ws.Range("A4").Value = $r4A
$r4C = @"
{'selector': <__main__.NamedFeatureSelector object at 0x7f3a6c394c70>, 'scaler': None, 'model__n_estimators': 10, 'model__estimator__solver': 'lbfgs', 'model__estimator__learning_rate_init': 0.01, 'model__estimator__hidden_layer_sizes': (10, 10, 10), 'model__estimator__alpha': 0.0001, 'model__estimator__activation': 'relu'}
"@
$ws.Range("C4").Value = $r4C
$r4E = @"
[1 0 1 1 1 1 0 1 0 1 0 1]
"@
$ws.Range("E4").Value = $r4E
$r4F = @"
[1 1 1 1 1 1 1 0 1 1 1 1]
"@
$ws.Range("F4").Value = $r4F
$ws.Range("B4").Value = 0.6095238095238095
$ws.Range("D4").Value = 0.7368421052631579
$ws.Range("G4").Value = 42
$ws.Range("H4").Value = 0.7530476190476191
$ws.Range("I4").Value = 0.02695810760902974
$ws.Range("J4").Value = 0.5393015873015873
$ws.Range("K4").Value = 0.07611533623491983

# Row 5
$r5A = @"
Pipeline(steps=[('scaler', None),
                ('selector',
                 SelectFromModel(estimator=LinearSVC(dual=False, penalty='l1',
                                                     random_state=42))),
                ('model',
                 BaggingClassifier(estimator=MLPClassifier(activation='tanh',
                                                           alpha=1e-05,
                                                           hidden_layer_sizes=(5,
                                                                               10,
                                                                               5),
                                                           max_iter=1000,
                                                           random_state=42,
                                                           solver='lbfgs'),
                                   n_estimators=5, random_state=42))])
"@
$ws.Range("A5").Value = $r5A
$r5C = @"
{'selector': SelectFromModel(estimator=LinearSVC(dual=False, penalty='l1', random_state=42)), 'scaler': None, 'model__n_estimators': 5, 'model__estimator__solver': 'lbfgs', 'model__estimator__learning_rate_init': 1, 'model__estimator__hidden_layer_sizes': (5, 10, 5), 'model__estimator__alpha': 1e-05, 'model__estimator__activation': 'tanh'}
"@
$ws.Range("C5").Value = $r5C
$r5E = @"
[1 1 0 0 0 0 1 0 1 1 1 1]
"@
$ws.Range("E5").Value = $r5E
$r5F = @"
[1 1 1 1 1 1 1 1 1 1 1 1]
"@
$ws.Range("F5").Value = $r5F
$ws.Range("B5").Value = 0.6285714285714286
$ws.Range("D5").Value = 0.7368421052631579
$ws.Range("G5").Value = 11
$ws.Range("H5").Value = 0.7590599876314162
$ws.Range("I5").Value = 0.03038426453614288
$ws.Range("J5").Value = 0.5257884972170687
$ws.Range("K5").Value = 0.07054333366411648

# Row 6
$r6A = @"
Pipeline(steps=[('scaler', None),
                ('selector',
                 SelectFromModel(estimator=LinearSVC(dual=False, penalty='l1',
                                                     random_state=42))),
                ('model',
                 BaggingClassifier(estimator=MLPClassifier(alpha=1,
                                                           hidden_layer_sizes=(10,),
                                                           learning_rate_init=0.01,
                                                           max_iter=1000,
                                                           random_state=42,
                                                           solver='sgd'),
                                   n_estimators=50, random_state=42))])
"@
$ws.Range("A6").Value = $r6A
$r6C = @"
{'selector': SelectFromModel(estimator=LinearSVC(dual=False, penalty='l1', random_state=42)), 'scaler': None, 'model__n_estimators': 50, 'model__estimator__solver': 'sgd', 'model__estimator__learning_rate_init': 0.01, 'model__estimator__hidden_layer_sizes': (10,), 'model__estimator__alpha': 1, 'model__estimator__activation': 'relu'}
"@
$ws.Range("C6").Value = $r6C
$r6E = @"
[1 1 1 1 0 0 0 0 1 1 0 0]
"@
$ws.Range("E6").Value = $r6E
$r6F = @"
[1 1 1 1 1 1 1 1 1 1 1 1]
"@
$ws.Range("F6").Value = $r6F
$ws.Range("B6").Value = 0.6190476190476191
$ws.Range("D6").Value = 0.6666666666666666
$ws.Range("G6").Value = 14
$ws.Range("H6").Value = 0.7742979242979242
$ws.Range("I6").Value = 0.02684926515624443
$ws.Range("J6").Value = 0.5628815628815628
$ws.Range("K6").Value = 0.07251142241405353
